$wb = $excel.ActiveWorkbook

# Map of row -> (old, new) for column F updates that apply to both
# "展览" and "全部类型" sheets.
$updates = @{
    2  = 1154
    3  = 871
    9  = 7808
    10 = 930
    11 = 452
    12 = 391
    13 = 161
    17 = 8024
    29 = 30
    31 = 1163
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}

# "全部类型" sheet has one additional change not present on "展览"
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F33").Value = 101
